# "Add files via upload" -- author populated the previously empty Sheet1
# with a small stock table (ticker / amount / name) and re-themed the
# accent1/accent5 colors.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate the stock table --------------------------------------------
# Written column-major (A down, then B, then C) so the shared-strings table
# is built up in the same order as in the authored workbook:
#   0 lvmhf, 1 nflx, 2 aapl, 3 amount, 4 ticker, 5 name, 6 Netflix, 7 Apple,
#   8 "louis vuitton moet hennessy"
$ws.Range("A2").Value = "lvmhf"
$ws.Range("A3").Value = "nflx"
$ws.Range("A4").Value = "aapl"
$ws.Range("B1").Value = "amount"
$ws.Range("A1").Value = "ticker"
$ws.Range("C1").Value = "name"
$ws.Range("C3").Value = "Netflix"
$ws.Range("C4").Value = "Apple"
$ws.Range("C2").Value = "louis vuitton moet hennessy"
$ws.Range("B2").Value = 5
$ws.Range("B3").Value = 4
$ws.Range("B4").Value = 7

# --- Column sizing & selection -------------------------------------------
# Column C ("name") was saved with a best-fit width in the authored file;
# AutoFit is the equivalent user action here.
$ws.Columns("C:C").AutoFit()

# The saved workbook has the populated range selected.
$null = $ws.Range("A1:C4").Select()

# --- Theme re-color: swap accent1 <-> accent5 -----------------------------
# Before: accent1=5B9BD5, accent5=4472C4
# After:  accent1=4472C4, accent5=5B9BD5
# ThemeColorScheme is indexed per MsoThemeColorSchemeIndex
# (1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1 ... 9 accent5, 10 accent6, 11 hlink, 12 folHlink)
# and RGB is an OLE color (0x00BBGGRR), so byte-swap each target RGB.
$tcs = $wb.Theme.ThemeColorScheme
$tcs.Colors(5).RGB = 0x00C47244   # -> accent1 srgbClr 4472C4
$tcs.Colors(9).RGB = 0x00D59B5B   # -> accent5 srgbClr 5B9BD5
